$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Small text edits on existing rows (no structural change)
# ---------------------------------------------------------------------------
$ws.Range("B5").Value  = 'The last time I saw my father......'
$ws.Range("B6").Value  = 'It must have been just before 6 PM. We had a huge argument, and I stormed out of the room.'

$ws.Range("B10").Value = ' <color=#00CC00>(To think the last time we met was in anger......)</color>'
$ws.Rows.Item(10).RowHeight = 34

$ws.Range("B11").Value = ' <color=#00CC00>(I wonder what Ming feels in his heart now. Is there regret?)</color>'

$ws.Range("B12").Value = 'At 6.15 PM, I saw signs of rain and headed toward the banquet hall.'
$ws.Range("J12").Value = 'appearAt'

$ws.Range("B13").Value = 'On the way, I ran into the newly arrived Kong, so we walked together.'
$ws.Range("B14").Value = 'I recall——it started raining just as you both arrived at the banquet hall?'

$ws.Range("B16").Value = ' <color=#00CC00>(Such a subtle difference in timing.)</color>'
$ws.Rows.Item(16).RowHeight = 34

$ws.Range("B17").Value = 'After the banquet started, I only left once——around 7.45 PM, I went to the backyard to urge Father to join us, but couldn’t find him.'

$ws.Range("B18").Value = 'Can you tell me more about your trip to the backyard?'
$ws.Rows.Item(18).RowHeight = 17

$ws.Range("J19").Value = 'appearAt'

# ---------------------------------------------------------------------------
# 2) Insert two new rows at position 21 (pushes old rows 21-24 to 23-26)
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

$ws.Rows.Item(21).RowHeight = 17
$ws.Range("A21").Value = "Dee"
$ws.Range("B21").Value = 'You went to parents’ bedroom?'
$ws.Range("C21").Value = "Dee-Thinking2"
$ws.Range("D21").Value = "DialogueVocal"
$ws.Range("E21").Value = "Question-Meeting"
$ws.Range("J21").Value = $null
$ws.Range("K21").Value = $null
$ws.Range("L21").Value = $null

$ws.Rows.Item(22).RowHeight = 51
$ws.Range("A22").Value = "Ming"
$ws.Range("B22").Value = 'Yes. At the time, Mother was probably bathing behind the screen——her damp clothes were hanging on the rack next to it.'
$ws.Range("C22").Value = "Ming-Regular"
$ws.Range("D22").Value = "DialogueVocal"
$ws.Range("E22").Value = "Question-Meeting"
$ws.Range("J22").Value = $null
$ws.Range("K22").Value = $null
$ws.Range("L22").Value = $null

# ---------------------------------------------------------------------------
# 3) Text edits on the rows that shifted down by 2 (now rows 23 & 24)
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = 'Not until 8.30 PM, when we all went together to look for him.'
$ws.Range("B24").Value = 'When you went to and from the backyard, did you pass through the frontyard?'

# ---------------------------------------------------------------------------
# 4) Sheet-level bookkeeping to match the final workbook state
# ---------------------------------------------------------------------------
$ws.Range("J22").Select()
